$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.2
$ws.Range("G3").Value = 26.7
$ws.Range("G4").Value = 21.3
$ws.Range("G5").Value = 24.9
$ws.Range("G6").Value = 19.100000000000001
$ws.Range("G7").Value = 22.4
$ws.Range("G8").Value = 18.600000000000001
$ws.Range("G10").Value = 25.6
$ws.Range("G11").Value = 28.1
$ws.Range("G12").Value = 18.600000000000001
$ws.Range("G13").Value = 22.2
$ws.Range("G14").Value = 25.3
$ws.Range("G15").Value = 29.5
$ws.Range("G16").Value = 17.399999999999999
$ws.Range("G17").Value = 20.9
$ws.Range("G18").Value = 24.7
$ws.Range("G20").Value = 27
$ws.Range("G21").Value = 30.2
$ws.Range("G22").Value = 26.5
$ws.Range("G23").Value = 52.6
$ws.Range("G24").Value = 18.899999999999999
$ws.Range("G25").Value = 36.1
$ws.Range("G26").Value = 23.9
$ws.Range("G27").Value = 47.5
$ws.Range("G28").Value = 26.9
$ws.Range("G29").Value = 51.8
$ws.Range("G30").Value = 21.1
$ws.Range("G31").Value = 46.7
$ws.Range("G32").Value = 18.100000000000001
$ws.Range("G33").Value = 33.200000000000003
$ws.Range("G34").Value = 18
$ws.Range("G35").Value = 27.8
$ws.Range("G36").Value = 41.9
$ws.Range("G37").Value = 24.4
$ws.Range("G38").Value = 53.5
$ws.Range("G39").Value = 24.9
$ws.Range("G41").Value = 23.2
$ws.Range("G42").Value = 48.2
$ws.Range("G43").Value = 23.8
$ws.Range("G44").Value = 49.6
$ws.Range("G45").Value = 25.3
$ws.Range("G48").Value = 26.6
$ws.Range("G49").Value = 50.5
$ws.Range("G50").Value = 19.8
$ws.Range("G51").Value = 45.9
$ws.Range("G52").Value = 17.7
$ws.Range("G53").Value = 36.4
$ws.Range("G54").Value = 17.600000000000001
$ws.Range("G55").Value = 30.5
$ws.Range("G56").Value = 18.2
$ws.Range("G57").Value = 38.4
$ws.Range("G58").Value = 20.3
$ws.Range("G59").Value = 36.700000000000003
$ws.Range("G60").Value = 21.9
$ws.Range("G61").Value = 48.8
$ws.Range("G62").Value = 24.9
$ws.Range("G63").Value = 24.5
$ws.Range("G64").Value = 53.8
$ws.Range("G65").Value = 26.4
$ws.Range("G66").Value = 51.9
$ws.Range("G67").Value = 17.100000000000001
$ws.Range("G68").Value = 27.5
$ws.Range("G69").Value = 18.100000000000001
$ws.Range("G70").Value = 43.3
$ws.Range("G71").Value = 25.7
$ws.Range("G72").Value = 49.5
$ws.Range("G73").Value = 19.399999999999999
$ws.Range("G74").Value = 20.2
$ws.Range("G75").Value = 46.2
$ws.Range("G76").Value = 19.100000000000001
$ws.Range("G77").Value = 34.4

$excel.ActiveWindow.ScrollRow = 39
